# Commit automatique via PowerShell
# Sorts the attendee table (rows 4-31, columns A-H) alphabetically by
# the "NOM" (last name) column (B), then renumbers column A back to a
# plain sequential 1..N index, matching the original numbering scheme.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 4
$lastRow = 31

$sortRange = $ws.Range("A$firstRow`:H$lastRow")
$sortKey = $ws.Range("B$firstRow`:B$lastRow")
$sortRange.Sort($sortKey, 1)

$count = $lastRow - $firstRow + 1
for ($i = 0; $i -lt $count; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 1).Value = $i + 1
}
